$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: 'data' ----------
$ws1 = $wb.Worksheets.Item("data")

# Copy AG1:AG67 formatting (incl. header style s=1) into the new AH1:AH67 column
$ws1.Range("AG1:AG67").Copy($ws1.Range("AH1:AH67"))

# New column header (date)
$ws1.Range("AH1").Value = "24. 8. 2021"

# New column AH values (rows 2-67)
$ws1.Range("AH2").Value = 0.09
$ws1.Range("AH3").Value = 0.18
$ws1.Range("AH4").Value = 0.73
$ws1.Range("AH5").Value = 0.08
$ws1.Range("AH6").Value = 0.14
$ws1.Range("AH7").Value = 0.78
$ws1.Range("AH8").Value = 0.04
$ws1.Range("AH9").Value = 0.12
$ws1.Range("AH10").Value = 0.84
$ws1.Range("AH11").Value = 0.09
$ws1.Range("AH12").Value = 0.22
$ws1.Range("AH13").Value = 0.6899999999999999
$ws1.Range("AH14").Value = 0.25
$ws1.Range("AH15").Value = 0.16
$ws1.Range("AH16").Value = 0.59
$ws1.Range("AH17").Value = 0.07000000000000001
$ws1.Range("AH18").Value = 0.14
$ws1.Range("AH19").Value = 0.79
$ws1.Range("AH20").Value = 0.05
$ws1.Range("AH21").Value = 0.11
$ws1.Range("AH22").Value = 0.84
$ws1.Range("AH23").Value = 0.2
$ws1.Range("AH24").Value = 0.17
$ws1.Range("AH25").Value = 0.63
$ws1.Range("AH26").Value = 0.11
$ws1.Range("AH27").Value = 0.21
$ws1.Range("AH28").Value = 0.68
$ws1.Range("AH29").Value = 0.15
$ws1.Range("AH30").Value = 0.26
$ws1.Range("AH31").Value = 0.59
$ws1.Range("AH32").Value = 0.06
$ws1.Range("AH33").Value = 0.14
$ws1.Range("AH34").Value = 0.8
$ws1.Range("AH35").Value = 0.03
$ws1.Range("AH36").Value = 0.06
$ws1.Range("AH37").Value = 0.91
$ws1.Range("AH38").Value = 0.16
$ws1.Range("AH39").Value = 0.19
$ws1.Range("AH40").Value = 0.65
$ws1.Range("AH41").Value = 0.07000000000000001
$ws1.Range("AH42").Value = 0.17
$ws1.Range("AH43").Value = 0.76
$ws1.Range("AH44").Value = 0.66
$ws1.Range("AH45").Value = 0.18
$ws1.Range("AH46").Value = 0.16
$ws1.Range("AH47").Value = 0.09
$ws1.Range("AH48").Value = 0.57
$ws1.Range("AH49").Value = 0.34
$ws1.Range("AH50").Value = 0.02
$ws1.Range("AH51").Value = 0.08
$ws1.Range("AH52").Value = 0.9
$ws1.Range("AH53").Value = 0.08
$ws1.Range("AH54").Value = 0.17
$ws1.Range("AH55").Value = 0.75
$ws1.Range("AH56").Value = 0.03
$ws1.Range("AH57").Value = 0.07000000000000001
$ws1.Range("AH58").Value = 0.9
$ws1.Range("AH59").Value = 0.08
$ws1.Range("AH60").Value = 0.2
$ws1.Range("AH61").Value = 0.72
$ws1.Range("AH62").Value = 0.05
$ws1.Range("AH63").Value = 0.08
$ws1.Range("AH64").Value = 0.87
$ws1.Range("AH65").Value = 0.06
$ws1.Range("AH66").Value = 0.06
$ws1.Range("AH67").Value = 0.88

# Corrections to existing AG column values (revised estimates)
$ws1.Range("AG3").Value = 0.17
$ws1.Range("AG4").Value = 0.73
$ws1.Range("AG6").Value = 0.13
$ws1.Range("AG7").Value = 0.8
$ws1.Range("AG39").Value = 0.19
$ws1.Range("AG40").Value = 0.67
$ws1.Range("AG45").Value = 0.16
$ws1.Range("AG46").Value = 0.16
$ws1.Range("AG48").Value = 0.55
$ws1.Range("AG49").Value = 0.35

# Update footer/title text with new update date
$ws1.Range("A68").Value = "Život během pandemie, Zasažení domácností, % respondentů celkově a ve skupinách, aktualizace 1. 9. 2021"

# ---------- Sheet 2: 'pocetR' ----------
$ws2 = $wb.Worksheets.Item("pocetR")

# Copy AF1:AF24 formatting (incl. header style s=2 and blank AF24) into the new AG1:AG24 column
$ws2.Range("AF1:AF24").Copy($ws2.Range("AG1:AG24"))

# New column header (date)
$ws2.Range("AG1").Value = "24. 8. 2021"

# New column AG values (rows 2-23); row 24 stays blank (footer row)
$ws2.Range("AG2").Value = 1620
$ws2.Range("AG3").Value = 768
$ws2.Range("AG4").Value = 132
$ws2.Range("AG5").Value = 505
$ws2.Range("AG6").Value = 215
$ws2.Range("AG7").Value = 726
$ws2.Range("AG8").Value = 124
$ws2.Range("AG9").Value = 105
$ws2.Range("AG10").Value = 665
$ws2.Range("AG11").Value = 759
$ws2.Range("AG12").Value = 546
$ws2.Range("AG13").Value = 315
$ws2.Range("AG14").Value = 453
$ws2.Range("AG15").Value = 1167
$ws2.Range("AG16").Value = 144
$ws2.Range("AG17").Value = 297
$ws2.Range("AG18").Value = 1179
$ws2.Range("AG19").Value = 285
$ws2.Range("AG20").Value = 90
$ws2.Range("AG21").Value = 272
$ws2.Range("AG22").Value = 149
$ws2.Range("AG23").Value = 77

# Update footer/title text with new update date
$ws2.Range("A24").Value = "Život během pandemie, Zasažení domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 9. 2021"
